$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-03-21"

# Update the row label for March to reflect the new "through" date
$ws.Range("A4").Value = "March (through 03-21)"

# Update March row (row 4) values for columns B, C, D, H, I
$ws.Range("B4").Value = 22
$ws.Range("C4").Value = 30
$ws.Range("D4").Value = 37
$ws.Range("H4").Value = 57
$ws.Range("I4").Value = 90

# Update Total row (row 5) values for columns B, C, D, H, I
$ws.Range("B5").Value = 59
$ws.Range("C5").Value = 117
$ws.Range("D5").Value = 168
$ws.Range("H5").Value = 399
$ws.Range("I5").Value = 390
